$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75, shifting rows 75-181 down to 76-182.
$ws.Range("A75").EntireRow.Insert()

# Populate the newly inserted row 75 with the new translation entry.
$ws.Range("A75").Value = "Initial and Final Surveillance Diagnosis"
$ws.Range("B75").Value = "TBT"

# Fix the bug with translation of "&": replace "&" with "and" in the
# (now shifted) row that used to read "Susceptible & Intermediate ..."
$ws.Range("A145").Value = "Susceptible and Intermediate are always combined in this visualisation of co-resistances."
